$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.239.66'
$ws.Range("E2").Value = '  -0.71%  '
$ws.Range("D3").Value = '1.902.27'
$ws.Range("E3").Value = '  +1.42%  '
$ws.Range("E4").Value = '  -0.48%  '
$ws.Range("B5").Value = 'XRP'
$ws.Range("C5").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.693'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +9.54%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '246.33'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.60%  '
$ws.Range("E7").Value = '  -0.36%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.92'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.32%  '
$ws.Range("E9").Value = '  +5.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '53.22'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +11.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0725'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0995'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D13").Value = '2.178.61'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '12.34'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.39%  '
$ws.Range("E15").Value = '  +3.28%  '
$ws.Range("D16").Value = '1.910.28'
$ws.Range("E16").Value = '  +1.57%  '
$ws.Range("E17").Value = '  +1.50%  '
$ws.Range("D18").Value = '35.239.57'
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.91%  '
$ws.Range("E20").Value = '  +2.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '241.18'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.77%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.61'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.19%  '
$ws.Range("E23").Value = '  +0.74%  '
$ws.Range("E24").Value = '  -0.38%  '
$ws.Range("E25").Value = '  +1.70%  '
$ws.Range("E26").Value = '  +15.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.87'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.84%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.52'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.07%  '
$ws.Range("E29").Value = '  +4.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '18.39'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.987'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.17'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.86%  '
$ws.Range("E34").Value = '  +1.29%  '
$ws.Range("E36").Value = '  +1.09%  '
$ws.Range("E37").Value = '  +0.79%  '
$ws.Range("E38").Value = '  -1.08%  '
$ws.Range("E39").Value = '  -0.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0686'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +16.30%  '
$ws.Range("E41").Value = '  +0.20%  '
$ws.Range("E42").Value = '  +3.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.29'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '90.79'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.09%  '
$ws.Range("D45").Value = '1.342.59'
$ws.Range("E45").Value = '  -0.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.43'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '47.19'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '12.59'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.41'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.47%  '
$ws.Range("E50").Value = '  +1.76%  '
$ws.Range("E51").Value = '  -2.11%  '
